# Apply odds updates for 2025-11-26 Betfair Back/Lay workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 1.35
$ws.Range("P2").Value = 1.35
$ws.Range("R2").Value = 1.2

# Row 3
$ws.Range("F3").Value = 1.34
$ws.Range("I3").Value = 12
$ws.Range("K3").Value = 5.9
$ws.Range("N3").Value = 5.8
$ws.Range("P3").Value = 2.64
$ws.Range("Q3").Value = 1.57
$ws.Range("R3").Value = 1.66
$ws.Range("S3").Value = 2.42
$ws.Range("U3").Value = 2.06
$ws.Range("W3").Value = 3.7
$ws.Range("AF3").Value = 8.800000000000001
$ws.Range("AM3").Value = 150
$ws.Range("AN3").Value = 4.7

# Row 4
$ws.Range("G4").Value = 5.5
$ws.Range("H4").Value = 1.73
$ws.Range("I4").Value = 1.74
$ws.Range("J4").Value = 4.1
$ws.Range("K4").Value = 4.3
$ws.Range("P4").Value = 2.34
$ws.Range("Q4").Value = 1.72
$ws.Range("S4").Value = 2.74
$ws.Range("T4").Value = 1.72
$ws.Range("V4").Value = 2.34
$ws.Range("Z4").Value = 11.5
$ws.Range("AC4").Value = 9.4
$ws.Range("AG4").Value = 20
$ws.Range("AH4").Value = 18
$ws.Range("AM4").Value = 85
$ws.Range("AO4").Value = 8

# Row 5
$ws.Range("I5").Value = 2.82
$ws.Range("L5").Value = 1.37
$ws.Range("N5").Value = 2.8
$ws.Range("P5").Value = 1.69
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.04

# Row 6
$ws.Range("L6").Value = 1.25
$ws.Range("X6").Value = 27
$ws.Range("Y6").Value = 24
$ws.Range("Z6").Value = 36
$ws.Range("AA6").Value = 65
$ws.Range("AB6").Value = 16.5
$ws.Range("AC6").Value = 11
$ws.Range("AD6").Value = 19.5
$ws.Range("AE6").Value = 44
$ws.Range("AF6").Value = 19.5
$ws.Range("AG6").Value = 13.5
$ws.Range("AH6").Value = 18
$ws.Range("AI6").Value = 48
$ws.Range("AJ6").Value = 32
$ws.Range("AK6").Value = 25
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 70
$ws.Range("AN6").Value = 13.5
$ws.Range("AO6").Value = 32

# Row 7
$ws.Range("F7").Value = 2.74
$ws.Range("G7").Value = 2.78
$ws.Range("I7").Value = 2.68
$ws.Range("L7").Value = 1.3
$ws.Range("V7").Value = 1.59
$ws.Range("W7").Value = 1.56
$ws.Range("X7").Value = 20
$ws.Range("Z7").Value = 19.5
$ws.Range("AA7").Value = 40
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 25
$ws.Range("AG7").Value = 12.5
$ws.Range("AI7").Value = 32
$ws.Range("AJ7").Value = 44
$ws.Range("AK7").Value = 26
$ws.Range("AL7").Value = 32
$ws.Range("AM7").Value = 65
$ws.Range("AO7").Value = 16.5

# Row 8
$ws.Range("F8").Value = 7.8
$ws.Range("G8").Value = 8.199999999999999
$ws.Range("H8").Value = 1.44
$ws.Range("I8").Value = 1.45
$ws.Range("J8").Value = 5.5
$ws.Range("L8").Value = 1.27
$ws.Range("N8").Value = 5.9
$ws.Range("P8").Value = 2.72
$ws.Range("R8").Value = 1.7
$ws.Range("V8").Value = 3.2
$ws.Range("W8").Value = 1.14
$ws.Range("X8").Value = 26
$ws.Range("Z8").Value = 10
$ws.Range("AC8").Value = 12.5
$ws.Range("AD8").Value = 10.5
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 75
$ws.Range("AG8").Value = 30
$ws.Range("AI8").Value = 27
$ws.Range("AJ8").Value = 260
$ws.Range("AK8").Value = 110
$ws.Range("AL8").Value = 85
$ws.Range("AM8").Value = 100
$ws.Range("AN8").Value = 95

# Row 9
$ws.Range("F9").Value = 2.32
$ws.Range("G9").Value = 2.36
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 3.35
$ws.Range("L9").Value = 1.34
$ws.Range("N9").Value = 4.6
$ws.Range("Q9").Value = 1.78
$ws.Range("V9").Value = 1.42
$ws.Range("W9").Value = 1.73
$ws.Range("X9").Value = 18
$ws.Range("Z9").Value = 26
$ws.Range("AA9").Value = 60
$ws.Range("AC9").Value = 8.199999999999999
$ws.Range("AD9").Value = 14
$ws.Range("AE9").Value = 34
$ws.Range("AF9").Value = 15.5
$ws.Range("AI9").Value = 40
$ws.Range("AL9").Value = 32
$ws.Range("AM9").Value = 70
$ws.Range("AO9").Value = 30

# Row 10
$ws.Range("G10").Value = 1.71
$ws.Range("L10").Value = 1.28
$ws.Range("U10").Value = 2.5
$ws.Range("V10").Value = 1.23
$ws.Range("W10").Value = 2.4
$ws.Range("X10").Value = 24
$ws.Range("Y10").Value = 25
$ws.Range("Z10").Value = 44
$ws.Range("AA10").Value = 120
$ws.Range("AC10").Value = 10
$ws.Range("AD10").Value = 20
$ws.Range("AE10").Value = 55
$ws.Range("AF10").Value = 12.5
$ws.Range("AG10").Value = 9.6
$ws.Range("AI10").Value = 55
$ws.Range("AJ10").Value = 17.5
$ws.Range("AK10").Value = 15.5
$ws.Range("AL10").Value = 26
$ws.Range("AM10").Value = 70
$ws.Range("AO10").Value = 44

# Row 11
$ws.Range("F11").Value = 1.39
$ws.Range("L11").Value = 1.27
$ws.Range("Q11").Value = 1.56
$ws.Range("R11").Value = 1.68
$ws.Range("U11").Value = 2.12
$ws.Range("V11").Value = 1.11
$ws.Range("W11").Value = 3.4
$ws.Range("X11").Value = 26
$ws.Range("Z11").Value = 85
$ws.Range("AD11").Value = 32
$ws.Range("AE11").Value = 120
$ws.Range("AF11").Value = 9.199999999999999
$ws.Range("AG11").Value = 10
$ws.Range("AI11").Value = 100
$ws.Range("AJ11").Value = 12
$ws.Range("AK11").Value = 13
$ws.Range("AL11").Value = 29
$ws.Range("AM11").Value = 110
$ws.Range("AO11").Value = 120

# Row 12
$ws.Range("G12").Value = 1.34
$ws.Range("H12").Value = 9.800000000000001
$ws.Range("L12").Value = 1.19
$ws.Range("N12").Value = 8.800000000000001
$ws.Range("O12").Value = 1.11
$ws.Range("P12").Value = 3.65
$ws.Range("R12").Value = 2.08
$ws.Range("S12").Value = 1.89
$ws.Range("T12").Value = 1.65
$ws.Range("V12").Value = 1.11
$ws.Range("W12").Value = 3.9
$ws.Range("X12").Value = 42
$ws.Range("Z12").Value = 110
$ws.Range("AA12").Value = 320
$ws.Range("AB12").Value = 16.5
$ws.Range("AC12").Value = 16
$ws.Range("AF12").Value = 12
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 12.5
$ws.Range("AK12").Value = 12.5
$ws.Range("AL12").Value = 24
$ws.Range("AN12").Value = 3.45

# Row 13
$ws.Range("K13").Value = 3.65
$ws.Range("L13").Value = 1.37
$ws.Range("O13").Value = 1.29
$ws.Range("V13").Value = 1.43
$ws.Range("W13").Value = 1.7
$ws.Range("X13").Value = 15
$ws.Range("Z13").Value = 23
$ws.Range("AA13").Value = 55
$ws.Range("AD13").Value = 14
$ws.Range("AE13").Value = 36
$ws.Range("AF13").Value = 15.5
$ws.Range("AG13").Value = 11
$ws.Range("AI13").Value = 48
$ws.Range("AK13").Value = 24
$ws.Range("AL13").Value = 40
$ws.Range("AM13").Value = 85
$ws.Range("AN13").Value = 17.5

# Row 14
$ws.Range("F14").Value = 2.04
$ws.Range("G14").Value = 2.16
$ws.Range("I14").Value = 4.5

# Row 15
$ws.Range("F15").Value = 1.77
$ws.Range("K15").Value = 4.3

